$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update defined names (rename + range adjustments)
$wb.Names.Item("actors").Delete()
$wb.Names.Add("Actors", "='Data Validation'!`$J`$5:`$J`$8")

$wb.Names.Item("categories").Delete()
$wb.Names.Add("Category", "='Data Validation'!`$B`$5:`$B`$8")

$wb.Names.Item("City").Delete()
$wb.Names.Add("City", "='Data Validation'!`$D`$5:`$D`$8")

$wb.Names.Item("food").Delete()
$wb.Names.Add("Food", "='Data Validation'!`$F`$5:`$F`$8")

$wb.Names.Item("fruit").Delete()
$wb.Names.Add("Fruit", "='Data Validation'!`$H`$5:`$H`$8")

# Add new content rows describing dependent dropdowns
$ws.Range("B48").Value = "creating dropdown of drop downs "
$ws.Range("E48").Value = "dependent drop downs)"
$ws.Range("H48").Value = "indirect(cell reference) is used "
$ws.Range("H49").Value = "to understand it once click on e51 and gothrough the value in it's data validation"
$ws.Range("C51").Value = "Fruit"
$ws.Range("E51").Value = "Mango"

# Update existing D40 validation to reference the renamed "Food" defined name
$ws.Range("D40").Validation.Delete()
$ws.Range("D40").Validation.Add(3, 1, 1, "Food")
$ws.Range("D40").Validation.ShowInput = $true
$ws.Range("D40").Validation.ShowError = $true

# Add new dropdown validations
$ws.Range("C51").Validation.Add(3, 1, 1, "Category")
$ws.Range("C51").Validation.ShowInput = $true
$ws.Range("C51").Validation.ShowError = $true

$ws.Range("E51").Validation.Add(3, 1, 1, '=INDIRECT($C$51)')
$ws.Range("E51").Validation.ShowInput = $true
$ws.Range("E51").Validation.ShowError = $true

# Update the view scroll position / selection
$excel.ActiveWindow.ScrollRow = 27
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H50").Select() | Out-Null
